$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Refreshed account-statement detail table (rows 16-78): worker identity
# (doc type/number/name), arrears period, arrears value and base salary,
# rebuilt from the updated source database.
$data = @(
    ,@('1047424130', 'OSNAIDER SALAS DIAZ', '2008', 25749, 877803)
    ,@('73184616', 'HENRY HERNANDEZ PRIMERA', '2101', 35112, 877803)
    ,@('73184616', 'HENRY HERNANDEZ PRIMERA', '2012', 35112, 877803)
    ,@('73184616', 'HENRY HERNANDEZ PRIMERA', '2011', 35112, 877803)
    ,@('73184616', 'HENRY HERNANDEZ PRIMERA', '2010', 35112, 877803)
    ,@('73184616', 'HENRY HERNANDEZ PRIMERA', '2009', 35112, 877803)
    ,@('73184616', 'HENRY HERNANDEZ PRIMERA', '2008', 35112, 877803)
    ,@('1047433490', 'ESTIVEN ALEXANDER MARIN CARMONA', '2102', 25749, 877803)
    ,@('1047433490', 'ESTIVEN ALEXANDER MARIN CARMONA', '2101', 35112, 877803)
    ,@('1047433490', 'ESTIVEN ALEXANDER MARIN CARMONA', '2012', 35112, 877803)
    ,@('1148434431', 'JORGE ALBERTO MASTRASCUSA VILLALOBOS', '2011', 35112, 877803)
    ,@('1047433490', 'ESTIVEN ALEXANDER MARIN CARMONA', '2010', 35112, 877803)
    ,@('1148434431', 'JORGE ALBERTO MASTRASCUSA VILLALOBOS', '2009', 35112, 877803)
    ,@('1047433490', 'ESTIVEN ALEXANDER MARIN CARMONA', '2008', 35112, 877803)
    ,@('1148434431', 'JORGE ALBERTO MASTRASCUSA VILLALOBOS', '2102', 25749, 877803)
    ,@('73115909', 'GUSTAVO FANEITE PEREZ', '2101', 35112, 877803)
    ,@('1148434431', 'JORGE ALBERTO MASTRASCUSA VILLALOBOS', '2012', 35112, 877803)
    ,@('1148434431', 'JORGE ALBERTO MASTRASCUSA VILLALOBOS', '2011', 35112, 877803)
    ,@('1148434431', 'JORGE ALBERTO MASTRASCUSA VILLALOBOS', '2010', 35112, 877803)
    ,@('73115909', 'GUSTAVO FANEITE PEREZ', '2009', 35112, 877803)
    ,@('1148434431', 'JORGE ALBERTO MASTRASCUSA VILLALOBOS', '2008', 35112, 877803)
    ,@('15681245', 'AQUILES ORTIZ LOPEZ', '2102', 25749, 877803)
    ,@('15681245', 'AQUILES ORTIZ LOPEZ', '2101', 35112, 877803)
    ,@('15681245', 'AQUILES ORTIZ LOPEZ', '2012', 35112, 877803)
    ,@('73115909', 'GUSTAVO FANEITE PEREZ', '2102', 25749, 877803)
    ,@('73115909', 'GUSTAVO FANEITE PEREZ', '2101', 35112, 877803)
    ,@('73115909', 'GUSTAVO FANEITE PEREZ', '2012', 35112, 877803)
    ,@('73115909', 'GUSTAVO FANEITE PEREZ', '2011', 35112, 877803)
    ,@('73115909', 'GUSTAVO FANEITE PEREZ', '2010', 35112, 877803)
    ,@('73115909', 'GUSTAVO FANEITE PEREZ', '2009', 35112, 877803)
    ,@('73115909', 'GUSTAVO FANEITE PEREZ', '2008', 35112, 877803)
    ,@('1047424130', 'OSNAIDER SALAS DIAZ', '2102', 25749, 877803)
    ,@('1047424130', 'OSNAIDER SALAS DIAZ', '2101', 35112, 877803)
    ,@('1047424130', 'OSNAIDER SALAS DIAZ', '2012', 35112, 877803)
    ,@('1047424130', 'OSNAIDER SALAS DIAZ', '2011', 35112, 877803)
    ,@('1047424130', 'OSNAIDER SALAS DIAZ', '2010', 35112, 877803)
    ,@('1047424130', 'OSNAIDER SALAS DIAZ', '2009', 35112, 877803)
    ,@('1047424130', 'OSNAIDER SALAS DIAZ', '2008', 35112, 877803)
    ,@('1048607162', 'ALVARO ENRIQUE MONTERO NARVAEZ', '2102', 25749, 877803)
    ,@('1048607162', 'ALVARO ENRIQUE MONTERO NARVAEZ', '2101', 35112, 877803)
    ,@('15617913', 'ANTONIO ENRIQUE MALO SOLAR', '2012', 35112, 877803)
    ,@('1048607162', 'ALVARO ENRIQUE MONTERO NARVAEZ', '2011', 35112, 877803)
    ,@('1048607162', 'ALVARO ENRIQUE MONTERO NARVAEZ', '2010', 35112, 877803)
    ,@('1048607162', 'ALVARO ENRIQUE MONTERO NARVAEZ', '2009', 35112, 877803)
    ,@('1048607162', 'ALVARO ENRIQUE MONTERO NARVAEZ', '2008', 35112, 877803)
    ,@('1131104611', 'ALEXANDER ARIAS MORELO', '2102', 25749, 877803)
    ,@('1131104893', 'BENJAMIN ENRIQUE HERRERA DIAZ', '2101', 35112, 877803)
    ,@('1131104893', 'BENJAMIN ENRIQUE HERRERA DIAZ', '2012', 35112, 877803)
    ,@('18810088', 'BENJAMIN ENRIQUE HERRERA SOLAR', '2102', 25749, 877803)
    ,@('18810088', 'BENJAMIN ENRIQUE HERRERA SOLAR', '2101', 35112, 877803)
    ,@('18810088', 'BENJAMIN ENRIQUE HERRERA SOLAR', '2012', 35112, 877803)
    ,@('92228882', 'EDWIN DE JESUS HERRERA SOLAR', '2102', 25749, 877803)
    ,@('92228882', 'EDWIN DE JESUS HERRERA SOLAR', '2101', 35112, 877803)
    ,@('92228882', 'EDWIN DE JESUS HERRERA SOLAR', '2012', 31601, 877803)
    ,@('1131104611', 'ALEXANDER ARIAS MORELO', '2102', 25749, 877803)
    ,@('1131104611', 'ALEXANDER ARIAS MORELO', '2101', 35112, 877803)
    ,@('1131104611', 'ALEXANDER ARIAS MORELO', '2012', 35112, 877803)
    ,@('15621880', 'JOSE FRANCISCO MENDOZA CANTERO', '2102', 25749, 877803)
    ,@('15621880', 'JOSE FRANCISCO MENDOZA CANTERO', '2101', 35112, 877803)
    ,@('15621880', 'JOSE FRANCISCO MENDOZA CANTERO', '2012', 35112, 877803)
    ,@('15617913', 'ANTONIO ENRIQUE MALO SOLAR', '2102', 25749, 877803)
    ,@('15617913', 'ANTONIO ENRIQUE MALO SOLAR', '2101', 35112, 877803)
    ,@('15617913', 'ANTONIO ENRIQUE MALO SOLAR', '2012', 35112, 877803)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 16 + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 3).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
}
